# Insert a new data row at row 384 (pushing the existing rows 384-482 down
# to 385-483), then populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(384).Insert()

$ws.Range("A384").Value = 10
$ws.Range("B384").Value = 'Vega Modelo de Temuco'
$ws.Range("C384").Value = 'La Araucanía'
$ws.Range("D384").Value = 44722
$ws.Range("E384").Value = 9
$ws.Range("F384").Value = 'Fruta'
$ws.Range("G384").Value = 100108
$ws.Range("H384").Value = 'Tropicales y subtropicales'
$ws.Range("I384").Value = 100108005
$ws.Range("J384").Value = 'Piña'
$ws.Range("K384").Value = 'Caramelo'
$ws.Range("L384").Value = 'Primera'
$ws.Range("M384").Value = 45
$ws.Range("N384").Value = 20000
$ws.Range("O384").Value = 20000
$ws.Range("P384").Value = 20000
$ws.Range("Q384").Value = '$/caja 14 unidades'
$ws.Range("R384").Value = 'Ecuador'
$ws.Range("S384").Value = 1429
$ws.Range("T384").Value = 14
